$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values per row (rows 2..18), reflecting the bugfix in the naive
# forecaster: each row's A/B/C/D values are shifted up from the row below
# (the oldest row is dropped), while the E (y_1_forecast) column is
# recomputed with corrected values.
$data = @(
  @{A=39765; B=2008; C=$null;                D=2009; E=3.407109591918855},
  @{A=40130; B=2009; C=-4.715480642250625;    D=2010; E=-2.725947775269033},
  @{A=40494; B=2010; C=6.130685532900904;     D=2011; E=5.993806847197725},
  @{A=40862; B=2011; C=8.703939237318981;     D=2012; E=6.924353497010971},
  @{A=41228; B=2012; C=2.688274587589135;     D=2013; E=3.737237443362851},
  @{A=41592; B=2013; C=0.9946838291217786;    D=2014; E=3.504647804006344},
  @{A=41957; B=2014; C=5.562499360312567;     D=2015; E=4.130094879572455},
  @{A=42321; B=2015; C=4.195080504802551;     D=2016; E=5.194458387461709},
  @{A=42689; B=2016; C=4.230623896992025;     D=2017; E=4.372655645302403},
  @{A=43053; B=2017; C=4.933871867981643;     D=2018; E=5.324897060120137},
  @{A=43418; B=2018; C=5.456119081407906;     D=2019; E=4.693063499664252},
  @{A=43783; B=2019; C=3.346849276607955;     D=2020; E=3.273620772016161},
  @{A=44159; B=2020; C=-9.2489161297999;      D=2021; E=-4.319815935184923},
  @{A=44525; B=2021; C=-1.287084480507283;    D=2022; E=1.850145327219943},
  @{A=44890; B=2022; C=1.494343500592232;     D=2023; E=2.552834403233084},
  @{A=45254; B=2023; C=-0.6982718287330991;   D=2024; E=-1.04392885455985},
  @{A=45618; B=2024; C=-0.4137309550271362;   D=2025; E=0.8860470190541037}
)

$rowIndex = 2
foreach ($row in $data) {
  $ws.Cells.Item($rowIndex, 1).Value = $row.A
  $ws.Cells.Item($rowIndex, 2).Value = $row.B
  if ($null -eq $row.C) {
    $ws.Cells.Item($rowIndex, 3).Value = $null
  } else {
    $ws.Cells.Item($rowIndex, 3).Value = $row.C
  }
  $ws.Cells.Item($rowIndex, 4).Value = $row.D
  $ws.Cells.Item($rowIndex, 5).Value = $row.E
  $rowIndex++
}

# Remove the now-obsolete last row (old row 19), shrinking the table to
# A1:E18.
$ws.Rows.Item(19).Delete()
